$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "summ28993959"
$wb.Worksheets.Item(2).Name = "summ29105748"
$wb.Worksheets.Item(3).Name = "summ29209169"
$wb.Worksheets.Item(4).Name = "summ29311532"
$wb.Worksheets.Item(5).Name = "summ29408860"
$wb.Worksheets.Item(6).Name = "summ29510731"
$wb.Worksheets.Item(7).Name = "summ29609552"
$wb.Worksheets.Item(8).Name = "summ29836543"
$wb.Worksheets.Item(9).Name = "summ29976107"
$wb.Worksheets.Item(10).Name = "summ30103428"
$wb.Worksheets.Item(11).Name = "summ30242880"
$wb.Worksheets.Item(12).Name = "summ30370636"
$wb.Worksheets.Item(13).Name = "summ30490494"
$wb.Worksheets.Item(14).Name = "summ30590266"
$wb.Worksheets.Item(15).Name = "summ30686057"
$wb.Worksheets.Item(16).Name = "summ30791776"
$wb.Worksheets.Item(17).Name = "summ30898540"
$wb.Worksheets.Item(18).Name = "summ30999328"
$wb.Worksheets.Item(19).Name = "summ31098079"
$wb.Worksheets.Item(20).Name = "summ31202357"
$wb.Worksheets.Item(21).Name = "summ31312346"
$wb.Worksheets.Item(22).Name = "summ31430393"
$wb.Worksheets.Item(23).Name = "summ31558039"
$wb.Worksheets.Item(24).Name = "summ31685060"
$wb.Worksheets.Item(25).Name = "summ31812234"
$wb.Worksheets.Item(26).Name = "summ31934030"
$wb.Worksheets.Item(27).Name = "summ32060862"
$wb.Worksheets.Item(28).Name = "summ32201230"
$wb.Worksheets.Item(29).Name = "summ32336921"
$wb.Worksheets.Item(30).Name = "summ32485879"
$wb.Worksheets.Item(31).Name = "summ32619716"
$wb.Worksheets.Item(32).Name = "summ32753536"
$wb.Worksheets.Item(33).Name = "summ32891605"
$wb.Worksheets.Item(34).Name = "summ33052852"
$wb.Worksheets.Item(35).Name = "summ33209475"
$wb.Worksheets.Item(36).Name = "summ33347051"
$wb.Worksheets.Item(37).Name = "summ33483439"
$wb.Worksheets.Item(38).Name = "summ33627006"
$wb.Worksheets.Item(39).Name = "summ33750497"
$wb.Worksheets.Item(40).Name = "summ33885273"
$wb.Worksheets.Item(41).Name = "summ34016493"
$wb.Worksheets.Item(42).Name = "summ34148043"
$wb.Worksheets.Item(43).Name = "summ34304960"
$wb.Worksheets.Item(44).Name = "summ34465136"
$wb.Worksheets.Item(45).Name = "summ34617639"
$wb.Worksheets.Item(46).Name = "summ34759505"
$wb.Worksheets.Item(47).Name = "summ34929930"
$wb.Worksheets.Item(48).Name = "summ35096358"
$wb.Worksheets.Item(49).Name = "summ35220581"
$wb.Worksheets.Item(50).Name = "summ35365470"
